# Apply updatePrice.js style updates to the inventory products sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: title change, price updates, checksum refresh
$ws.Range("B2").Value = "گردنبند طلا زنانه"
$ws.Range("I2").Value = 12073284
$ws.Range("J2").Value = 13901400
$ws.Range("U2").Value = "3d0df742b5b4"

# Row 3: no longer buybox winner -> price drops, buybox price cleared, flag flips
$ws.Range("F3").Value = 0
$ws.Range("I3").Value = 12281444
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = "خیر"
$ws.Range("U3").Value = "e6723b848790"

# Row 4
$ws.Range("F4").Value = 0
$ws.Range("I4").Value = 12990789
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = "خیر"
$ws.Range("U4").Value = "f645dc1c02a6"

# Row 5
$ws.Range("F5").Value = 0
$ws.Range("I5").Value = 13226170
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = "خیر"
$ws.Range("U5").Value = "ee6ebc510dbb"

# Row 6
$ws.Range("F6").Value = 0
$ws.Range("I6").Value = 11516055
$ws.Range("J6").ClearContents()
$ws.Range("K6").Value = "خیر"
$ws.Range("U6").Value = "fac4c114bf4d"

# Row 7
$ws.Range("F7").Value = 0
$ws.Range("I7").Value = 14076424
$ws.Range("J7").ClearContents()
$ws.Range("K7").Value = "خیر"
$ws.Range("U7").Value = "1ffc27d247d8"

# Row 8
$ws.Range("F8").Value = 0
$ws.Range("I8").Value = 8497734
$ws.Range("J8").ClearContents()
$ws.Range("K8").Value = "خیر"
$ws.Range("U8").Value = "dc5486bd1e10"

# Row 9
$ws.Range("F9").Value = 0
$ws.Range("I9").Value = 22200000
$ws.Range("J9").ClearContents()
$ws.Range("K9").Value = "خیر"
$ws.Range("U9").Value = "889469343ae2"
